$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B40: change from text "3" to numeric 3
$ws.Range("B40").Value = 3

# Add new row 41 with annotation data
$ws.Range("A41").Value = "Sunsi Wu"
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "2"
$ws.Range("C41").Value = "however"
$ws.Range("D41").Value = "ACK"
$ws.Range("E41").Value = "RES"
$ws.Range("F41").Value = "42b1e2ab-785d-481e-b197-1cf6913a8b3e"
$ws.Range("G41").Value = "SJQO7UJCW_annotated.xlsx"
$ws.Range("H41").Value = "However, our main point of the paper is to demonstrate the effectiveness of proposed method against our baseline model shown in Table 1 and 2."
